$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Planilha2")

# Update the value in Planilha1!A2
$ws1.Range("A2").Value = "joao123998"

# Selections / active sheet
# Set the selection on Planilha2 first (without leaving it as the active sheet)
$ws2.Select()
$ws2.Range("D6").Select()

# Then activate Planilha1 and select A2 there, leaving it as the active tab
$ws1.Select()
$ws1.Range("A2").Select()
